$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.773812157043043
$ws.Range("C2").Value = -0.6716102162172155
$ws.Range("D2").Value = 0.6737542003636109

$ws.Range("B3").Value = 0.7734895351277546
$ws.Range("C3").Value = -0.8074185313179805
$ws.Range("D3").Value = 0.7769351374832036

$ws.Range("B4").Value = -0.6960891120842669
$ws.Range("C4").Value = 0.7959025124123852
$ws.Range("D4").Value = -0.7452268809369937

$ws.Range("B5").Value = 0.811071421775988
$ws.Range("C5").Value = -0.7137521536241178
$ws.Range("D5").Value = 0.621472011254339

$ws.Range("B6").Value = -0.7685360354389394
$ws.Range("C6").Value = -0.7610006807898221
$ws.Range("D6").Value = -0.6825160850545469

$ws.Range("B7").Value = -0.6796641061219179
$ws.Range("C7").Value = -0.6130205986050971
$ws.Range("D7").Value = -0.6068112116314853

$ws.Range("B8").Value = 0.7203592411142159
$ws.Range("C8").Value = 0.5948621725398051
$ws.Range("D8").Value = -0.5957808337143961

$ws.Range("B9").Value = 0.6815708191199694
$ws.Range("C9").Value = -0.7920251443836751
$ws.Range("D9").Value = 0.7304928071655132
